# Correction of data scale (division by 1000) of acoustic data from visual experiment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Call-frequency")

# Divide all values in columns D, E, F (rows 2-40) by 1000
$rng = $ws.Range("D2:F40")
foreach ($cell in $rng.Cells) {
    $cell.Value2 = $cell.Value2 / 1000
}

# Best-fit width on column D (as observed in the saved file, ~9.57 characters)
$ws.Columns.Item(4).ColumnWidth = 8.74

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("H11").Select()
